$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3753
$ws.Range("I86").Value = 3841.25
$ws.Range("J86").Value = 3488.25
$ws.Range("K86").Value = 3841.25
$ws.Range("L86").Value = 3488.25
$ws.Range("M86").Value = -2718.25
$ws.Range("N86").Value = -5734.25
$ws.Range("H89").Value = 3753
$ws.Range("I89").Value = 3841.25
$ws.Range("J89").Value = 3488.25
$ws.Range("K89").Value = 19206.25
$ws.Range("L89").Value = 17441.25
$ws.Range("M89").Value = -13590.25
$ws.Range("N89").Value = -28673.25
$ws.Range("H98").Value = 666.3333
$ws.Range("I98").Value = 666.3333
$ws.Range("K98").Value = 666.3333
$ws.Range("M98").Value = 831.6667
$ws.Range("H112").Value = 2464.2307
$ws.Range("I112").Value = 1290
$ws.Range("J112").Value = 2677.7273
$ws.Range("K112").Value = 3870
$ws.Range("L112").Value = 8033.1819
$ws.Range("M112").Value = -2762
$ws.Range("N112").Value = -10249.1819
$ws.Range("H116").Value = 4104.6
$ws.Range("I116").Value = 4506.375
$ws.Range("K116").Value = 4506.375
$ws.Range("M116").Value = -1064.375
$ws.Range("H122").Value = 666.3333
$ws.Range("I122").Value = 666.3333
$ws.Range("K122").Value = 1998.9999
$ws.Range("M122").Value = 451.0001
$ws.Range("H137").Value = 15153004
$ws.Range("I137").Value = 25642122
$ws.Range("J137").Value = 2055.3333
$ws.Range("K137").Value = 76926366
$ws.Range("L137").Value = 6165.999899999999
$ws.Range("M137").Value = -76923816
$ws.Range("N137").Value = -11265.9999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6381.275
$ws.Range("I32").Value = 4856.7295
$ws.Range("K32").Value = 4856.7295
$ws.Range("M32").Value = -4569.7295
$ws.Range("H63").Value = 2100
$ws.Range("I63").Value = 2100
$ws.Range("K63").Value = 2100
$ws.Range("M63").Value = -1414
$ws.Range("H66").Value = 2100
$ws.Range("I66").Value = 2100
$ws.Range("K66").Value = 10500
$ws.Range("M66").Value = -7068
$ws.Range("H88").Value = 787.44446
$ws.Range("I88").Value = 399
$ws.Range("J88").Value = 1273
$ws.Range("K88").Value = 399
$ws.Range("L88").Value = 1273
$ws.Range("M88").Value = 7
$ws.Range("N88").Value = -2085
$ws.Range("H91").Value = 787.44446
$ws.Range("I91").Value = 399
$ws.Range("J91").Value = 1273
$ws.Range("K91").Value = 399
$ws.Range("L91").Value = 1273
$ws.Range("M91").Value = 1005
$ws.Range("N91").Value = -4081
$ws.Range("H102").Value = 5066.387
$ws.Range("I102").Value = 3803.9546
$ws.Range("K102").Value = 3803.9546
$ws.Range("M102").Value = -2181.9546
$ws.Range("H122").Value = 6835.4165
$ws.Range("I122").Value = 6902.8
$ws.Range("J122").Value = 6498.5
$ws.Range("K122").Value = 20708.4
$ws.Range("L122").Value = 19495.5
$ws.Range("M122").Value = -18258.4
$ws.Range("N122").Value = -24395.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 181.625
$ws.Range("I5").Value = 88
$ws.Range("J5").Value = 462.5
$ws.Range("K5").Value = 88
$ws.Range("L5").Value = 462.5
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = -688.5
$ws.Range("H35").Value = 47536.75
$ws.Range("I35").Value = 44999
$ws.Range("J35").Value = 50074.5
$ws.Range("K35").Value = 44999
$ws.Range("L35").Value = 50074.5
$ws.Range("M35").Value = -44689
$ws.Range("N35").Value = -50694.5
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H86").Value = 12376446
$ws.Range("I86").Value = 33106.562
$ws.Range("J86").Value = 30330394
$ws.Range("K86").Value = 33106.562
$ws.Range("L86").Value = 30330394
$ws.Range("M86").Value = -31983.562
$ws.Range("N86").Value = -30332640
$ws.Range("H89").Value = 12376446
$ws.Range("I89").Value = 33106.562
$ws.Range("J89").Value = 30330394
$ws.Range("K89").Value = 165532.81
$ws.Range("L89").Value = 151651970
$ws.Range("M89").Value = -159916.81
$ws.Range("N89").Value = -151663202
$ws.Range("H99").Value = 1933.3334
$ws.Range("I99").Value = 1933.3334
$ws.Range("K99").Value = 1933.3334
$ws.Range("M99").Value = -435.3334
$ws.Range("H107").Value = 2991.5
$ws.Range("I107").Value = 3057.6
$ws.Range("K107").Value = 3057.6
$ws.Range("M107").Value = -1137.6

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1956.9231
$ws.Range("I31").Value = 1679.4736
$ws.Range("J31").Value = 2710
$ws.Range("K31").Value = 1679.4736
$ws.Range("L31").Value = 2710
$ws.Range("M31").Value = -1384.4736
$ws.Range("N31").Value = -3300
$ws.Range("H34").Value = 1956.9231
$ws.Range("I34").Value = 1679.4736
$ws.Range("J34").Value = 2710
$ws.Range("K34").Value = 1679.4736
$ws.Range("L34").Value = 2710
$ws.Range("M34").Value = -1477.4736
$ws.Range("N34").Value = -3114
$ws.Range("H62").Value = 2433
$ws.Range("I62").Value = 2399.75
$ws.Range("J62").Value = 2499.5
$ws.Range("K62").Value = 2399.75
$ws.Range("L62").Value = 2499.5
$ws.Range("M62").Value = -1775.75
$ws.Range("N62").Value = -3747.5
$ws.Range("H65").Value = 2433
$ws.Range("I65").Value = 2399.75
$ws.Range("J65").Value = 2499.5
$ws.Range("K65").Value = 11998.75
$ws.Range("L65").Value = 12497.5
$ws.Range("M65").Value = -8878.75
$ws.Range("N65").Value = -18737.5
$ws.Range("H96").Value = 15228
$ws.Range("J96").Value = 15228
$ws.Range("L96").Value = 15228
$ws.Range("N96").Value = -20720
$ws.Range("H112").Value = 54995
$ws.Range("J112").Value = 54995
$ws.Range("L112").Value = 54995
$ws.Range("N112").Value = -57949
$ws.Range("H141").Value = 122475.625
$ws.Range("J141").Value = 122475.625
$ws.Range("L141").Value = 122475.625
$ws.Range("N141").Value = -132835.625

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 88667.164
$ws.Range("I34").Value = 312.75
$ws.Range("J34").Value = 132844.38
$ws.Range("K34").Value = 938.25
$ws.Range("L34").Value = 398533.14
$ws.Range("M34").Value = -854.25
$ws.Range("N34").Value = -398701.14
$ws.Range("H55").Value = 328
$ws.Range("I55").Value = 328
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 984
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -807
$ws.Range("N55").ClearContents()
$ws.Range("H113").Value = 2812.5715
$ws.Range("I113").Value = 2629
$ws.Range("J113").Value = 2886
$ws.Range("K113").Value = 7887
$ws.Range("L113").Value = 8658
$ws.Range("M113").Value = -5717
$ws.Range("N113").Value = -12998

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 2989
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 2989
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 2989
$ws.Range("N6").Value = -3215
$ws.Range("M6").ClearContents()
$ws.Range("H16").Value = 2989
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2989
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2989
$ws.Range("N16").Value = -3489
$ws.Range("M16").ClearContents()
$ws.Range("H80").Value = 2492.2856
$ws.Range("J80").Value = 2492.2856
$ws.Range("L80").Value = 2492.2856
$ws.Range("N80").Value = -4488.2856
$ws.Range("H83").Value = 2492.2856
$ws.Range("J83").Value = 2492.2856
$ws.Range("L83").Value = 12461.428
$ws.Range("N83").Value = -22445.428

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 722.2308
$ws.Range("I16").Value = 722.2308
$ws.Range("K16").Value = 722.2308
$ws.Range("M16").Value = -552.2308
$ws.Range("H22").Value = 100000840
$ws.Range("I22").Value = 731
$ws.Range("K22").Value = 731
$ws.Range("M22").Value = -436
$ws.Range("H27").Value = 100000840
$ws.Range("I27").Value = 731
$ws.Range("K27").Value = 731
$ws.Range("M27").Value = -624
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H82").Value = 1578.6666
$ws.Range("J82").Value = 1833
$ws.Range("L82").Value = 1833
$ws.Range("N82").Value = -2555
$ws.Range("H85").Value = 1578.6666
$ws.Range("J85").Value = 1833
$ws.Range("L85").Value = 1833
$ws.Range("N85").Value = -4329
$ws.Range("H94").Value = 68388.164
$ws.Range("J94").Value = 68388.164
$ws.Range("L94").Value = 68388.164
$ws.Range("N94").Value = -69740.164
$ws.Range("H132").Value = 8053.8887
$ws.Range("J132").Value = 50000
$ws.Range("L132").Value = 150000
$ws.Range("N132").Value = -155060
$ws.Range("H136").Value = 2198.3845
$ws.Range("J136").Value = 2433.3333
$ws.Range("L136").Value = 7299.999899999999
$ws.Range("N136").Value = -12399.9999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H81").Value = 8070.387
$ws.Range("I81").Value = 4009.8948
$ws.Range("K81").Value = 8019.7896
$ws.Range("M81").Value = -6958.7896
$ws.Range("H84").Value = 8070.387
$ws.Range("I84").Value = 4009.8948
$ws.Range("K84").Value = 40098.948
$ws.Range("M84").Value = -34794.948
